$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oligos")

# --- Insert a new "Tag" column at the front of the Oligos sheet ---
$ws.Columns.Item(1).Insert()

# Header
$ws.Range("A1").Value = "Tag"

# Tag values (base oligo id, derived from the "Name" column in B)
$ws.Range("A2").Value = "o2"
$ws.Range("A3").Value = "o3"
$ws.Range("A4").Value = "o4"

# --- Column widths (character units); shifted one column to the right ---
# (target widths, in the saved-file's units, are 9.8/32.19/9.98/17.47/15.37/9.98/8.06/9.4 -
#  the values below are the ColumnWidth inputs that round-trip closest to those)
$ws.Columns.Item(1).ColumnWidth = 9.0
$ws.Columns.Item(2).ColumnWidth = 9.0
$ws.Columns.Item(3).ColumnWidth = 31.333333
$ws.Columns.Item(4).ColumnWidth = 31.333333
$ws.Columns.Item(5).ColumnWidth = 9.166667
$ws.Columns.Item(6).ColumnWidth = 16.666667
$ws.Columns.Item(7).ColumnWidth = 14.5
$ws.Columns.Item(8).ColumnWidth = 9.166667
$ws.Columns.Item(9).ColumnWidth = 7.166667
$ws.Columns.Item(10).ColumnWidth = 8.5

# --- Update the Print_Area named range for Oligos (shifted one column right) ---
$n = $ws.Names.Item("Oligos!Print_Area")
$n.RefersTo = '=Oligos!$I$94:$J$98'

# --- Move the selection like the edited file did ---
$ws.Range("A5").Select()
